# Commit: "Fruta / hortaliza, semanal"
#
# Inserts 4 new weekly price-report rows for "Comercializadora del Agro de
# Limarí - Naranja" on Sheet1, above the existing row 210. Excel pushes the
# previously-existing rows 210:289 down to 214:293 (the sheet's used range
# grows from A1:T289 to A1:T293), and the new rows carry a later report date
# (2021-09-16, serial 44455) for varieties "Lane Late" and "Navel Late".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 210; this shifts the old
# 210:289 block down to 214:293 and grows the sheet dimension accordingly.
$ws.Rows("210:213").Insert()

# Data for the 4 newly inserted rows (same 20-column layout as every other
# record row: A..T = Mercado ID .. Kg/unidad).
$newRows = @(
    @{ Row = 210; Variedad = 'Lane Late';  Calidad = 'Primera'; Volumen = 20; Min = 110000; Max = 120000; Prom = 115000; KgKg = 288 },
    @{ Row = 211; Variedad = 'Lane Late';  Calidad = 'Segunda'; Volumen = 20; Min =  90000; Max = 100000; Prom =  95000; KgKg = 238 },
    @{ Row = 212; Variedad = 'Navel Late'; Calidad = 'Primera'; Volumen = 20; Min = 110000; Max = 120000; Prom = 115000; KgKg = 288 },
    @{ Row = 213; Variedad = 'Navel Late'; Calidad = 'Segunda'; Volumen = 20; Min =  90000; Max = 100000; Prom =  95000; KgKg = 238 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 2                                     # A Mercado ID
    $ws.Cells.Item($row, 2).Value  = 'Comercializadora del Agro de Limarí' # B Mercado
    $ws.Cells.Item($row, 3).Value  = 'Coquimbo'                            # C Región
    $ws.Cells.Item($row, 4).Value  = 44455                                 # D Fecha (2021-09-16)
    $ws.Cells.Item($row, 5).Value  = 4                                     # E Codreg
    $ws.Cells.Item($row, 6).Value  = 'Fruta'                               # F Tipo
    $ws.Cells.Item($row, 7).Value  = 100102                                # G Producto ID
    $ws.Cells.Item($row, 8).Value  = 'Cítricos'                            # H Producto
    $ws.Cells.Item($row, 9).Value  = 100102005                             # I Categoría ID
    $ws.Cells.Item($row, 10).Value = 'Naranja'                             # J Categoría
    $ws.Cells.Item($row, 11).Value = $r.Variedad                           # K Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad                            # L Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen                            # M Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min                                # N Precio mínimo
    $ws.Cells.Item($row, 15).Value = $r.Max                                # O Precio máximo
    $ws.Cells.Item($row, 16).Value = $r.Prom                               # P Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = '$/bins (400 kilos)'                  # Q Unidad de comercialización
    $ws.Cells.Item($row, 18).Value = 'Provincia de Limarí'                 # R Origen
    $ws.Cells.Item($row, 19).Value = $r.KgKg                               # S Precio $/Kg
    $ws.Cells.Item($row, 20).Value = 400                                   # T Kg / unidad
}
